$wb = $excel.ActiveWorkbook

# --- Sheet1: "Tactic upload" ---
$ws1 = $wb.Worksheets.Item(1)

# Add two new empty rows (B2, B3) formatted as dates (builtin numFmtId 14).
# Format B2 directly, then copy/paste-special (formats only) onto B3 so both
# cells share a single cellXfs entry instead of each write minting its own.
$ws1.Cells.Item(2, 2).NumberFormat = "mm-dd-yy"
$ws1.Cells.Item(2, 2).Copy()
$ws1.Cells.Item(3, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Sheet2: "Channels for Reference" ---
$ws2 = $wb.Worksheets.Item(2)

$values = @(
  "Circuit_HR_Pages",
  "Circuit_News",
  "Circuit_Microsite",
  "Wordpress_Microsite",
  "Ask_Vote_Answer",
  "Double_Dutch",
  "Physical_Poster",
  "Digital_Sign",
  "Email",
  "Inside_Blue",
  "MeetUp",
  "My_Intel_App",
  "Intel_Newsroom",
  "SharePoint",
  "Webcast",
  "Enterprise_Wiki",
  "Twitter",
  "LinkedIn",
  "Facebook",
  "Instagram"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Record sheet2's saved selection (A20), then switch back to sheet1 so it
# remains the active/selected tab, matching the workbook's saved state.
$ws2.Activate()
[void]$ws2.Range("A20").Select()

$ws1.Activate()
[void]$ws1.Range("G12").Select()
